$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the Animals / Activities / 2020 (COVID-theme) word lists ---
# Values are written in the same order the shared-string table records them
# (first pass fills column A top-to-bottom, then column B, then column C,
#  and finally the newly inserted row-3 animal is typed last).

$ws.Range("A2").Value = "MOKNEY"
$ws.Range("A4").Value = "CAT"
$ws.Range("A5").Value = "DOG"
$ws.Range("A6").Value = "BIRD"
$ws.Range("A7").Value = "HIPPO"
$ws.Range("A8").Value = "COW"
$ws.Range("A9").Value = "HORSE"
$ws.Range("A10").Value = "WHALE"
$ws.Range("A11").Value = "TURTLE"
$ws.Range("A12").Value = "FISH"
$ws.Range("A13").Value = "GOAT"
$ws.Range("A14").Value = "RABBIT"
$ws.Range("A15").Value = "SNAKE"

$ws.Range("B2").Value = "SKI"
$ws.Range("B3").Value = "RUN"
$ws.Range("B4").Value = "JUMP"
$ws.Range("B5").Value = "SING"
$ws.Range("B6").Value = "DANCE"
$ws.Range("B7").Value = "WRITE"
$ws.Range("B9").Value = "DRIVE"
$ws.Range("B8").Value = "BIKE"
$ws.Range("B10").Value = "PARTY"
$ws.Range("B11").Value = "EAT"
$ws.Range("B12").Value = "DRINK"
$ws.Range("B13").Value = "CLEAN"
$ws.Range("B14").Value = "STUDY"
$ws.Range("B15").Value = "DIE"

$ws.Range("C2").Value = "QUARANTINE"
$ws.Range("C3").Value = "COVID"
$ws.Range("C4").Value = "CORONA"
$ws.Range("C5").Value = "BLACK LIVES MATTER"
$ws.Range("C6").Value = "WWIII"

# A new animal row was inserted at row 3 - typed in last.
$ws.Range("A3").Value = "RHINO"

# --- Formatting: extend the existing (row 8) cell style down through A15 ---
$ws.Range("A9:A15").Font.Bold = $true

# --- Selection / scroll position housekeeping ---
[void]$ws.Range("A3").Select()
